# Updates cryptos list values (prices, 1h volume %) and swaps Aave/Hedera row order
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'60.488.81"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.27%  '
$ws.Range('D3').Value = "'2.637.51"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.88%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = "'567.80"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.89%  '
$ws.Range('D6').Value = "'146.48"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.50%  '
$ws.Range('D7').Value = "'0.996"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('E8').Value = '  +4.02%  '
$ws.Range('D9').Value = "'2.662.65"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.66%  '
$ws.Range('E10').Value = '  +0.13%  '
$ws.Range('E11').Value = '  +4.78%  '
$ws.Range('E12').Value = '  +6.77%  '
$ws.Range('E13').Value = '  +2.84%  '
$ws.Range('D14').Value = "'3.109.62"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.19%  '
$ws.Range('D15').Value = "'60.468.24"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.33%  '
$ws.Range('D16').Value = "'22.07"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.05%  '
$ws.Range('E17').Value = '  +4.32%  '
$ws.Range('D18').Value = "'2.657.64"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.66%  '
$ws.Range('E19').Value = '  +2.90%  '
$ws.Range('D20').Value = "'343.02"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.35%  '
$ws.Range('D21').Value = "'10.44"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.93%  '
$ws.Range('D22').Value = "'6.38"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.51%  '
$ws.Range('E24').Value = '  -1.54%  '
$ws.Range('E25').Value = '  +4.54%  '
$ws.Range('E26').Value = '  +2.44%  '
$ws.Range('D27').Value = "'0.991"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('E28').Value = '  +4.47%  '
$ws.Range('D29').Value = "'0.0₃0809"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +9.91%  '
$ws.Range('D30').Value = "'0.997"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.15%  '
$ws.Range('D31').Value = "'1.72"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.23%  '
$ws.Range('E32').Value = '  +4.82%  '
$ws.Range('D33').Value = "'159.54"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.25%  '
$ws.Range('E34').Value = '  +1.82%  '
$ws.Range('E35').Value = '  +5.07%  '
$ws.Range('D36').Value = "'0.904"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +7.95%  '
$ws.Range('E37').Value = '  +4.73%  '
$ws.Range('E38').Value = '  +8.71%  '
$ws.Range('E39').Value = '  +6.72%  '
$ws.Range('D40').Value = "'37.45"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.91%  '
$ws.Range('D41').Value = "'300.97"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.97%  '
$ws.Range('E42').Value = '  +1.35%  '
$ws.Range('D43').Value = "'0.995"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.58%  '
$ws.Range('D44').Value = "'0.0987"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.74%  '
$ws.Range('D45').Value = "'0.605"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.70%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'129.35"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +15.60%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').Value = "'0.0546"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.84%  '
$ws.Range('E48').Value = '  +1.30%  '
$ws.Range('E49').Value = '  +0.04%  '
$ws.Range('E50').Value = '  +4.36%  '
$ws.Range('D51').Value = "'4.66"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.01%  '
